$wb = $excel.ActiveWorkbook

# --- PLO sheet: content edits -------------------------------------------------
$wsPLO = $wb.Worksheets.Item("PLO")

# Existing PLO row (row 5) description in English gains an "(Optional)" suffix.
$wsPLO.Range("C5").Value = "Able to… (Optional)"

# Fill in the blank "insert row" of Table1 (A2:B3) with the Thai label for the
# Year column - this converts the table's virtual insert row into real data.
$wsPLO.Range("B3").Value = "ปีของหลักสูตร"

# The dropdown for the curriculum-type cell changes its allowed values.
$wsPLO.Range("A3").Validation.Formula1 = '"regular,inter"'

# --- Selections on the other sheets (set while they're not the active sheet,
#     so they don't steal tabSelected / activeTab) ------------------------------
$wsPO = $wb.Worksheets.Item("PO")
$wsPO.Range("C3").Select()

# --- Make PLO the active sheet & set its final selection -----------------------
$wsPLO.Activate()
$wsPLO.Range("B12").Select()
